$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Summary paragraph -- reorder / extend the technology list.
#   "...including: Bootstrap, HTML5, CSS3, JavaScript, AutoCAD, Jira, Revit,
#    SQL Queries & SQL Server 2008-2016. "
#   becomes
#   "...including: T-SQL Queries & SQL Server 2008-2016, Archibus Web
#    Central, Bootstrap, HTML5, CSS3, JavaScript, AutoCAD, Revit, Jira."
#   (the trailing formatted space run is dropped and the _GoBack bookmark
#   ends up right after the new final "Jira.")
# ---------------------------------------------------------------------------

$oldSummary = "including: Bootstrap, HTML5, CSS3, JavaScript, AutoCAD, Jira, Revit, SQL Queries & SQL Server 2008-2016."
$newSummary = "including: T-SQL Queries & SQL Server 2008-2016, Archibus Web Central, Bootstrap, HTML5, CSS3, JavaScript, AutoCAD, Revit, Jira."

$d.Content.Find.Execute($oldSummary, $true, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2)

# Re-anchor the _GoBack bookmark right after the new "...Revit, Jira." text.
$bmRng = $d.Content
$bmRng.Find.Execute("Revit, Jira.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jiraEnd = $bmRng.End
$bmRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Drop the trailing formatted space run that used to close the paragraph.
$trailingSpace = $d.Range($jiraEnd, $jiraEnd + 1)
$trailingSpace.Delete()

# ---------------------------------------------------------------------------
# Edit 2: Job title heading -- reorder the role list.
#   "Front-End Web Developer, CAD Designer, IWMS Administrator "
#   becomes
#   "IWMS Administrator, CAD Designer, Frontend Web Developer"
# ---------------------------------------------------------------------------

$oldTitle = "Front-End Web Developer, CAD Designer, IWMS Administrator "
$newTitle = "IWMS Administrator, CAD Designer, Frontend Web Developer"
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)

# ---------------------------------------------------------------------------
# Edit 3: Bullet point -- drop "jQuery/" before "JavaScript".
# ---------------------------------------------------------------------------

$oldBullet = "• Developed over 25 Responsive SpaceView Bootstrap webpages with 100+ design and functionality updates using HTML5, CSS3, and jQuery/JavaScript."
$newBullet = "• Developed over 25 Responsive SpaceView Bootstrap webpages with 100+ design and functionality updates using HTML5, CSS3, and JavaScript."
$d.Content.Find.Execute($oldBullet, $true, $false, $false, $false, $false, $true, 1, $false, $newBullet, 2)
